# Insert 3 new price-report rows ("Clemenuless" variety) at the top of the
# Mandarina table (rows 130-132), pushing the existing rows 130-212 down to
# 133-215. The sheet's used range grows from A1:T212 to A1:T215.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 130..212 down by 3 (creates 3 blank rows at 130:132, carrying
# the date-format style already applied to column D down through the sheet).
$ws.Rows("130:132").Insert()

# --- Row 130: Clemenuless / Especial -----------------------------------
$ws.Range("A130").Value = 4
$ws.Range("B130").Value = 'Feria Lagunitas de Puerto Montt'
$ws.Range("C130").Value = 'Los Lagos'
$ws.Range("D130").Value = 44719
$ws.Range("E130").Value = 10
$ws.Range("F130").Value = 'Fruta'
$ws.Range("G130").Value = 100102
$ws.Range("H130").Value = 'Cítricos'
$ws.Range("I130").Value = 100102004
$ws.Range("J130").Value = 'Mandarina'
$ws.Range("K130").Value = 'Clemenuless'
$ws.Range("L130").Value = 'Especial'
$ws.Range("M130").Value = 300
$ws.Range("N130").Value = 12500
$ws.Range("O130").Value = 12500
$ws.Range("P130").Value = 12500
$ws.Range("Q130").Value = '$/bandeja 10 kilos'
$ws.Range("R130").Value = 'Provincia de Limarí'
$ws.Range("S130").Value = 1250
$ws.Range("T130").Value = 10

# --- Row 131: Clemenuless / Primera -------------------------------------
$ws.Range("A131").Value = 4
$ws.Range("B131").Value = 'Feria Lagunitas de Puerto Montt'
$ws.Range("C131").Value = 'Los Lagos'
$ws.Range("D131").Value = 44719
$ws.Range("E131").Value = 10
$ws.Range("F131").Value = 'Fruta'
$ws.Range("G131").Value = 100102
$ws.Range("H131").Value = 'Cítricos'
$ws.Range("I131").Value = 100102004
$ws.Range("J131").Value = 'Mandarina'
$ws.Range("K131").Value = 'Clemenuless'
$ws.Range("L131").Value = 'Primera'
$ws.Range("M131").Value = 300
$ws.Range("N131").Value = 10000
$ws.Range("O131").Value = 10000
$ws.Range("P131").Value = 10000
$ws.Range("Q131").Value = '$/bandeja 10 kilos'
$ws.Range("R131").Value = 'Provincia de Limarí'
$ws.Range("S131").Value = 1000
$ws.Range("T131").Value = 10

# --- Row 132: Clemenuless / Segunda -------------------------------------
$ws.Range("A132").Value = 4
$ws.Range("B132").Value = 'Feria Lagunitas de Puerto Montt'
$ws.Range("C132").Value = 'Los Lagos'
$ws.Range("D132").Value = 44719
$ws.Range("E132").Value = 10
$ws.Range("F132").Value = 'Fruta'
$ws.Range("G132").Value = 100102
$ws.Range("H132").Value = 'Cítricos'
$ws.Range("I132").Value = 100102004
$ws.Range("J132").Value = 'Mandarina'
$ws.Range("K132").Value = 'Clemenuless'
$ws.Range("L132").Value = 'Segunda'
$ws.Range("M132").Value = 300
$ws.Range("N132").Value = 8500
$ws.Range("O132").Value = 8500
$ws.Range("P132").Value = 8500
$ws.Range("Q132").Value = '$/bandeja 10 kilos'
$ws.Range("R132").Value = 'Provincia de Limarí'
$ws.Range("S132").Value = 850
$ws.Range("T132").Value = 10
